$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch the affected range first so the text number-format sticks,
# then restore normal style after all values are written, to avoid
# leaving a stray "quote-prefixed"/custom-format style on the cells.
$priceVolRange = $ws.Range("D2:E47")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "301.00"
$ws.Range("E2").Value = "-0.84%"
$ws.Range("D3").Value = "31.39"
$ws.Range("D4").Value = "5.095"
$ws.Range("E4").Value = "-3.15%"
$ws.Range("D5").Value = "0.07372"
$ws.Range("E5").Value = "-1.89%"
$ws.Range("D6").Value = "2.381"
$ws.Range("E6").Value = "57.40%"
$ws.Range("D7").Value = "7.962"
$ws.Range("E7").Value = "1.02%"
$ws.Range("D8").Value = "3.790"
$ws.Range("E8").Value = "-0.74%"
$ws.Range("D9").Value = "0.9166"
$ws.Range("E9").Value = "-0.55%"
$ws.Range("D10").Value = "0.1713"
$ws.Range("E10").Value = "0.47%"
$ws.Range("D11").Value = "0.07534"
$ws.Range("E11").Value = "-4.70%"
$ws.Range("D12").Value = "0.08053"
$ws.Range("E12").Value = "0.36%"
$ws.Range("D13").Value = "0.03015"
$ws.Range("E13").Value = "-1.02%"
$ws.Range("E14").Value = "0.24%"
$ws.Range("D15").Value = "0.001497"
$ws.Range("E15").Value = "0.49%"
$ws.Range("D16").Value = "0.006175"
$ws.Range("E16").Value = "-2.65%"
$ws.Range("D17").Value = "3.455"
$ws.Range("E17").Value = "-0.20%"
$ws.Range("D18").Value = "2.225"
$ws.Range("E18").Value = "-0.35%"
$ws.Range("D19").Value = "0.3296"
$ws.Range("E19").Value = "-0.25%"
$ws.Range("D20").Value = "0.1337"
$ws.Range("E20").Value = "-0.41%"
$ws.Range("D21").Value = "4.647"
$ws.Range("E21").Value = "3.48%"
$ws.Range("D22").Value = "0.04644"
$ws.Range("E22").Value = "0.75%"
$ws.Range("D23").Value = "0.1566"
$ws.Range("E23").Value = "-3.22%"
$ws.Range("D24").Value = "0.001226"
$ws.Range("E24").Value = "0.77%"
$ws.Range("D25").Value = "0.004489"
$ws.Range("E25").Value = "0.92%"
$ws.Range("E26").Value = "-7.04%"
$ws.Range("E27").Value = "50.46%"
$ws.Range("D39").Value = "0.01739"
$ws.Range("E39").Value = "0.73%"
$ws.Range("E40").Value = "0.84%"
$ws.Range("D41").Value = "0.007219"
$ws.Range("E41").Value = "3.78%"
$ws.Range("D42").Value = "0.1348"
$ws.Range("E42").Value = "-0.08%"
$ws.Range("E43").Value = "1.47%"
$ws.Range("D44").Value = "0.01074"
$ws.Range("E44").Value = "-16.83%"
$ws.Range("D45").Value = "0.00006299"
$ws.Range("E45").Value = "2.21%"
$ws.Range("E46").Value = "-33.27%"
$ws.Range("D47").Value = "0.8085"
$ws.Range("E47").Value = "-56.66%"

$priceVolRange.Style = "Normal"

